# MatrizGestionInteresados.xlsx - apply the "Añadidos apartados a la memoria
# y modificada la matriz de estrategias" edits to the Hoja1 worksheet:
#  - widen column E (more room for the updated "Estrategia" text)
#  - grow several row heights to fit the newly-added/expanded text
#  - move the active selection to G3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E (Estrategia) to match columns A:B's width (22.21875
# characters). The host's column-width setter quantises to the nearest
# 1/6-character increment, so feed it the input that lands closest to the
# true target (the nearest reachable stored width is 22.1666...).
$ws.Columns.Item(5).ColumnWidth = 21.33

# Row height adjustments (explicit custom heights).
$ws.Rows.Item(1).RowHeight = 52.8
$ws.Rows.Item(2).RowHeight = 34.8
$ws.Rows.Item(12).RowHeight = 57
$ws.Rows.Item(21).RowHeight = 40.2
$ws.Rows.Item(22).RowHeight = 46.2
$ws.Rows.Item(29).RowHeight = 51

# Move the selection/active cell to G3.
$ws.Range("G3").Select()
